# Update "想去人数" (number of people interested) values in column F
# for both the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet,
# which contain duplicated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 6636
    $ws.Range("F6").Value = 2028
    $ws.Range("F7").Value = 1541
    $ws.Range("F9").Value = 1011
    $ws.Range("F12").Value = 5635
}
